$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Components")

# --- New row 5: 1uF 0603 capacitor ---
$ws.Range("L5").Value = "1uF Caps"
$ws.Range("K5").Value = "C0603C105K8PACTU"
$ws.Hyperlinks.Add($ws.Range("M5"), "https://www.digikey.com.au/product-detail/en/kemet/C0603C105K8PACTU/399-3118-1-ND/551623")
$ws.Range("M5").Style = "Hyperlink"

# --- New row 6: 0.1uF 0603 capacitor ---
$ws.Range("L6").Value = "0.1uF Caps"
$ws.Hyperlinks.Add($ws.Range("M6"), "https://www.digikey.com.au/product-detail/en/kemet/C0603C104Z3VACTU/399-1100-1-ND/411375")
$ws.Range("M6").Style = "Hyperlink"
$ws.Range("K6").Value = "C0603C104Z3VACTU"

# --- Stray space entries ---
$ws.Range("P3").Value = " "
$ws.Range("P4").Value = " "
$ws.Range("P5").Value = " "
$ws.Range("P6").Value = " "

# --- Final selection ---
$ws.Range("P7:R7").Select() | Out-Null
